$d = $word.ActiveDocument

# The BTec logo inline pictures (in the headers) were renamed image1.jpg -> image2.jpg,
# and the Pearson Edexcel logo inline pictures (in the footers) were renamed
# image2.png -> image1.png. Walk every section's headers/footers and rename each
# inline picture based on which logo it is (identified by its stable AlternativeText /
# picture description, since header/footer index<->part mapping can vary).

for ($s = 1; $s -le $d.Sections.Count; $s++) {
  $sec = $d.Sections.Item($s)

  for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
      $shapes = $hdr.Range.InlineShapes
      for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
          $shp.Name = "image2.jpg"
        }
      }
    }

    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
      $shapes = $ftr.Range.InlineShapes
      for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
          $shp.Name = "image1.png"
        }
      }
    }
  }
}
